# This script reproduces a periodic data refresh of the "上海-漫展信息"
# (Shanghai convention-info) workbook: mostly small increments to the
# "want to go" counter (column F) across all four sheets, plus one
# outdated listing on sheet 3 ("本地生活") that expired and was removed,
# shifting the rows below it up by one.

$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 1: 展览 (Exhibitions) - refresh "want to go" counts (col F)
# ============================================================
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(4, 6).Value = 450
$ws1.Cells.Item(5, 6).Value = 1861
$ws1.Cells.Item(7, 6).Value = 7431
$ws1.Cells.Item(8, 6).Value = 204
$ws1.Cells.Item(10, 6).Value = 200
$ws1.Cells.Item(11, 6).Value = 1686
$ws1.Cells.Item(12, 6).Value = 1442
$ws1.Cells.Item(13, 6).Value = 1270
$ws1.Cells.Item(14, 6).Value = 134
$ws1.Cells.Item(15, 6).Value = 134
$ws1.Cells.Item(16, 6).Value = 3454
$ws1.Cells.Item(17, 6).Value = 5817
$ws1.Cells.Item(18, 6).Value = 5817
$ws1.Cells.Item(20, 6).Value = 607
$ws1.Cells.Item(21, 6).Value = 968
$ws1.Cells.Item(22, 6).Value = 1190
$ws1.Cells.Item(23, 6).Value = 353
$ws1.Cells.Item(24, 6).Value = 5841
$ws1.Cells.Item(25, 6).Value = 332
$ws1.Cells.Item(27, 6).Value = 51
$ws1.Cells.Item(28, 6).Value = 4003
$ws1.Cells.Item(29, 6).Value = 215
$ws1.Cells.Item(30, 6).Value = 670
$ws1.Cells.Item(31, 6).Value = 1847
$ws1.Cells.Item(32, 6).Value = 1125
$ws1.Cells.Item(35, 6).Value = 160
$ws1.Cells.Item(36, 6).Value = 107
$ws1.Cells.Item(37, 6).Value = 305
$ws1.Cells.Item(38, 6).Value = 1108
$ws1.Cells.Item(39, 6).Value = 476
$ws1.Cells.Item(40, 6).Value = 1817
$ws1.Cells.Item(41, 6).Value = 79
$ws1.Cells.Item(42, 6).Value = 356
$ws1.Cells.Item(43, 6).Value = 132
$ws1.Cells.Item(44, 6).Value = 1027
$ws1.Cells.Item(47, 6).Value = 57
$ws1.Cells.Item(50, 6).Value = 140

# ============================================================
# Sheet 2: 演出 (Performances) - refresh "want to go" counts (col F)
# ============================================================
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(4, 6).Value = 977
$ws2.Cells.Item(5, 6).Value = 116
$ws2.Cells.Item(7, 6).Value = 18
$ws2.Cells.Item(8, 6).Value = 12
$ws2.Cells.Item(11, 6).Value = 325
$ws2.Cells.Item(12, 6).Value = 384
$ws2.Cells.Item(13, 6).Value = 22
$ws2.Cells.Item(14, 6).Value = 192
$ws2.Cells.Item(15, 6).Value = 100
$ws2.Cells.Item(16, 6).Value = 7
$ws2.Cells.Item(18, 6).Value = 337
$ws2.Cells.Item(19, 6).Value = 132
$ws2.Cells.Item(20, 6).Value = 141
$ws2.Cells.Item(22, 6).Value = 226
$ws2.Cells.Item(23, 6).Value = 71
$ws2.Cells.Item(26, 6).Value = 26
$ws2.Cells.Item(27, 6).Value = 251
$ws2.Cells.Item(34, 6).Value = 3

# ============================================================
# Sheet 3: 本地生活 (Local life)
#   - row 8 ("怪兽8号 niko and ... 集章之旅") expired and is removed;
#     rows 9-15 shift up to become rows 8-14.
#   - column A is a plain sequential index (0-based row number) that is
#     NOT supposed to travel with the shifted content, so it is restored
#     to 7..13 for the new rows 8..14 after the shift.
#   - "want to go" counts (col F) are refreshed for every surviving row.
# ============================================================
$ws3 = $wb.Worksheets.Item(3)

# F4 / F5 / F7 are unaffected by the row shift (rows above the deletion).
$ws3.Cells.Item(4, 6).Value = 3321
$ws3.Cells.Item(5, 6).Value = 436
$ws3.Cells.Item(7, 6).Value = 1521

# Remove the expired row 8; Excel shifts rows 9-15 up to 8-14 automatically,
# carrying B:I content with them and updating the sheet dimension.
$ws3.Rows.Item(8).Delete()

# Column A must stay a fixed sequential index per row (A8=7 ... A14=13),
# independent of which listing now occupies that row.
for ($r = 8; $r -le 14; $r++) {
    $ws3.Cells.Item($r, 1).Value = $r - 1
}

# Refresh "want to go" counts (col F) for the shifted rows (now 8-14).
$ws3.Cells.Item(8, 6).Value = 449
$ws3.Cells.Item(9, 6).Value = 3017
$ws3.Cells.Item(10, 6).Value = 387
$ws3.Cells.Item(11, 6).Value = 771
$ws3.Cells.Item(12, 6).Value = 936
$ws3.Cells.Item(13, 6).Value = 1000
$ws3.Cells.Item(14, 6).Value = 1444

# ============================================================
# Sheet 4: 全部类型 (All types) - refresh "want to go" counts (col F)
# ============================================================
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 436
$ws4.Cells.Item(4, 6).Value = 450
$ws4.Cells.Item(5, 6).Value = 449
$ws4.Cells.Item(6, 6).Value = 3017
$ws4.Cells.Item(8, 6).Value = 7431
$ws4.Cells.Item(9, 6).Value = 204
$ws4.Cells.Item(10, 6).Value = 18
$ws4.Cells.Item(11, 6).Value = 771
$ws4.Cells.Item(13, 6).Value = 200
$ws4.Cells.Item(14, 6).Value = 1442
$ws4.Cells.Item(15, 6).Value = 1000
$ws4.Cells.Item(17, 6).Value = 134
$ws4.Cells.Item(18, 6).Value = 134
$ws4.Cells.Item(19, 6).Value = 3454
$ws4.Cells.Item(20, 6).Value = 325
$ws4.Cells.Item(21, 6).Value = 5817
$ws4.Cells.Item(22, 6).Value = 22
$ws4.Cells.Item(23, 6).Value = 607
$ws4.Cells.Item(24, 6).Value = 968
$ws4.Cells.Item(25, 6).Value = 1190
$ws4.Cells.Item(26, 6).Value = 353
$ws4.Cells.Item(27, 6).Value = 5841
$ws4.Cells.Item(28, 6).Value = 332
$ws4.Cells.Item(29, 6).Value = 4003
$ws4.Cells.Item(30, 6).Value = 670
$ws4.Cells.Item(31, 6).Value = 337
$ws4.Cells.Item(32, 6).Value = 1847
$ws4.Cells.Item(33, 6).Value = 1125
$ws4.Cells.Item(35, 6).Value = 132
$ws4.Cells.Item(36, 6).Value = 160
$ws4.Cells.Item(37, 6).Value = 107
$ws4.Cells.Item(38, 6).Value = 305
$ws4.Cells.Item(39, 6).Value = 1108
$ws4.Cells.Item(40, 6).Value = 1817
$ws4.Cells.Item(41, 6).Value = 79
$ws4.Cells.Item(42, 6).Value = 356
$ws4.Cells.Item(43, 6).Value = 132
$ws4.Cells.Item(44, 6).Value = 1027
$ws4.Cells.Item(47, 6).Value = 251
$ws4.Cells.Item(50, 6).Value = 140
